$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Propagate the existing cell style (date number format / border / alignment)
# from the last populated row (A343) down through the new rows (A344:A357),
# exactly like dragging the fill handle down in Excel.
$ws.Range("A343").Copy($ws.Range("A344:A357"))

# New daily data rows (r=344..357), covering 2021-08-10 .. 2021-08-23
# (date serials 44418..44431), matching the update described in the commit
# message ("aggiornamento al 23 agosto 2021").
$data = @(
    @(44418, 0, 4, 26.39218791237794),
    @(44419, 0, 4, 26.39218791237794),
    @(44420, 1, 4, 26.39218791237794),
    @(44421, 3, 6, 39.5882818685669),
    @(44422, 3, 9, 59.38242280285036),
    @(44423, 1, 8, 52.78437582475588),
    @(44424, 1, 9, 59.38242280285036),
    @(44425, 2, 11, 72.57851675903932),
    @(44426, 0, 11, 72.57851675903932),
    @(44427, 0, 10, 65.98046978094484),
    @(44428, 5, 12, 79.1765637371338),
    @(44429, 0, 9, 59.38242280285036),
    @(44430, 5, 13, 85.77461071522829),
    @(44431, 0, 12, 79.1765637371338)
)

$startRow = 344
for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}
